$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '44.191.85'
$ws.Range("E2").Value = '  +5.45%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.267.78'
$ws.Range("E3").Value = '  +2.46%  '

$ws.Range("E4").Value = '  +0.23%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '231.00'
$ws.Range("E5").Value = '  +0.22%  '

$ws.Range("E6").Value = '  +2.88%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '63.93'
$ws.Range("E7").Value = '  +5.31%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.441'
$ws.Range("E9").Value = '  +9.88%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.104'
$ws.Range("E10").Value = '  +15.43%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '56.58'
$ws.Range("E11").Value = '  -1.08%  '

$ws.Range("B12").Value = 'Avalanche'
$ws.Range("C12").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '26.15'
$ws.Range("E12").Value = '  +18.29%  '

$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.106'
$ws.Range("E13").Value = '  +2.69%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.605.10'
$ws.Range("E14").Value = '  +2.54%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '15.77'
$ws.Range("E15").Value = '  +2.19%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.09'
$ws.Range("E16").Value = '  +9.37%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.835'
$ws.Range("E17").Value = '  +4.83%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.254.49'
$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '44.085.88'
$ws.Range("E19").Value = '  +5.39%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0000103'
$ws.Range("E20").Value = '  +11.34%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '73.68'
$ws.Range("E21").Value = '  +2.36%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  -0.23%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '256.09'
$ws.Range("E23").Value = '  +5.58%  '

$ws.Range("E24").Value = '  -0.05%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.44'
$ws.Range("E25").Value = '  +3.90%  '

$ws.Range("E26").Value = '  -6.22%  '

$ws.Range("B27").Value = 'WEMIXToken'
$ws.Range("C27").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '3.37'
$ws.Range("E27").Value = '  +27.91%  '

$ws.Range("B28").Value = 'Cosmos'
$ws.Range("C28").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.06'
$ws.Range("E28").Value = '  +4.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '172.12'
$ws.Range("E29").Value = '  +1.89%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.92'
$ws.Range("E30").Value = '  +5.79%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.138'
$ws.Range("E31").Value = '  -1.23%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.40'
$ws.Range("E32").Value = '  -2.43%  '

$ws.Range("E33").Value = '  +2.75%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0684'
$ws.Range("E34").Value = '  +5.25%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.75'
$ws.Range("E35").Value = '  +3.22%  '

$ws.Range("E36").Value = '  -1.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.88'
$ws.Range("E37").Value = '  +9.68%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.75'
$ws.Range("E38").Value = '  +6.70%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.33'
$ws.Range("E39").Value = '  +0.37%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0258'
$ws.Range("E40").Value = '  +6.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.37'
$ws.Range("E42").Value = '  -1.33%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '17.43'
$ws.Range("E43").Value = '  +8.16%  '

$ws.Range("B44").Value = 'TerraClassic'
$ws.Range("C44").Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.000216'
$ws.Range("E44").Value = '  -7.81%  '

$ws.Range("B45").Value = 'Cronos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0960'
$ws.Range("E45").Value = '  +0.83%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '97.85'
$ws.Range("E46").Value = '  +1.18%  '

$ws.Range("E47").Value = '  -0.19%  '

$ws.Range("B48").Value = 'FTXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '4.39'
$ws.Range("E48").Value = '  +1.83%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.448.62'
$ws.Range("E49").Value = '  -0.15%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.30'
$ws.Range("E50").Value = '  +3.83%  '

$ws.Range("B51").Value = 'Celestia'
$ws.Range("C51").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '9.87'
$ws.Range("E51").Value = '  +16.01%  '

